# TC02_Bento_Filter_PRStatus-NotReported.xlsx - "Fixed Bento 80 Test scripts"
#
# Appends " order By <col> ASC LIMIT 100" clauses to the three Cypher queries
# stored in column B (CasesTab / SamplesTab / FilesTab rows), and grows the
# row heights of rows 2 & 3 to fit the extra line of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2 (CasesTab query): append a new trailing line with ORDER BY / LIMIT ---
$b2Old = $ws.Range("B2").Value2
$b2Suffix = "`n order By ss.study_subject_id ASC LIMIT 100 "
if ($b2Old.EndsWith($b2Suffix) -eq $false) {
    $ws.Range("B2").Value = $b2Old + $b2Suffix
}

# --- B3 (SamplesTab query): append a new trailing line with ORDER BY / LIMIT ---
$b3Old = $ws.Range("B3").Value2
$b3Suffix = "`n order By samp.sample_id ASC LIMIT 100"
if ($b3Old.EndsWith($b3Suffix) -eq $false) {
    $ws.Range("B3").Value = $b3Old + $b3Suffix
}

# --- B4 (FilesTab query): replace the old "order by f.file_name" tail line ---
$b4Old = $ws.Range("B4").Value2
$b4OldTail = "    order by f.file_name"
$b4NewTail = "     order By f.file_name ASC LIMIT 100"
if ($b4Old.EndsWith($b4OldTail)) {
    $b4New = $b4Old.Substring(0, $b4Old.Length - $b4OldTail.Length) + $b4NewTail
    $ws.Range("B4").Value = $b4New
}

# --- Row heights grew (345.6 -> 360) to accommodate the extra query line ---
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 360

# --- Scroll the view down so row 3 is at the top (best-effort; mirrors the
#     workbook being scrolled/re-viewed after the edit) ---
$win = $excel.ActiveWindow
if ($win -ne $null) {
    $win.ScrollRow = 3
    $win.ScrollColumn = 1
}
